# Append the latest EUR->ARS quote as a new row at the bottom of the sheet.
# The existing data is stored as plain text (dates like "2025-10-22" and
# times like "21:20:38" are literal strings, not date/time serials), so we
# temporarily force a text number format before writing the values and then
# restore the default "Normal" style afterwards so no stray style index is
# left on the new cells (matching the rest of the sheet, which has none).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
if ($newRow -lt 2) { $newRow = 2 }

$rng = $ws.Range("A$newRow`:B$newRow")
$rng.NumberFormat = "@"

$ws.Range("A$newRow").Value = "2025-10-22"
$ws.Range("B$newRow").Value = "21:20:38"
$ws.Range("C$newRow").Value = "1.00 EUR = 1,842.4025"

$rng.Style = "Normal"
